$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 24-34 are brand new data rows appended to the list; give column A
# the same bold/bordered/centered "index number" look used by rows 2-23
# (copy formatting only, values are set explicitly afterwards).
$ws.Range("A2").Copy()
$ws.Range("A24:A34").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 2
$ws.Cells.Item(2, 2).Value = "NSE:ANUP"
$ws.Cells.Item(2, 3).Value = "NSE:20MICRONS"
$ws.Cells.Item(2, 4).ClearContents()
$ws.Cells.Item(2, 5).ClearContents()
$ws.Cells.Item(2, 6).Value = "NSE:HDFCLIFE"

# Row 3
$ws.Cells.Item(3, 2).Value = "NSE:CAPLIPOINT"
$ws.Cells.Item(3, 3).Value = "NSE:AHLUCONT"
$ws.Cells.Item(3, 5).ClearContents()
$ws.Cells.Item(3, 6).ClearContents()

# Row 4
$ws.Cells.Item(4, 2).Value = "NSE:HDFCLIFE"
$ws.Cells.Item(4, 3).Value = "NSE:AMRUTANJAN"

# Row 5
$ws.Cells.Item(5, 2).Value = "NSE:HESTERBIO"
$ws.Cells.Item(5, 3).Value = "NSE:ANANDRATHI"

# Row 6
$ws.Cells.Item(6, 2).Value = "NSE:HONASA"
$ws.Cells.Item(6, 3).Value = "NSE:BAJAJHIND"

# Row 7
$ws.Cells.Item(7, 2).Value = "NSE:INDIANB"
$ws.Cells.Item(7, 3).Value = "NSE:CDSL"

# Row 8
$ws.Cells.Item(8, 2).Value = "NSE:JMFINANCIL"
$ws.Cells.Item(8, 3).Value = "NSE:CLEAN"

# Row 9
$ws.Cells.Item(9, 2).Value = "NSE:KALYANKJIL"
$ws.Cells.Item(9, 3).Value = "NSE:CONSOFINVT"

# Row 10
$ws.Cells.Item(10, 2).Value = "NSE:NIFTYETF"
$ws.Cells.Item(10, 3).Value = "NSE:DATAPATTNS"

# Row 11
$ws.Cells.Item(11, 2).Value = "NSE:NITINSPIN"
$ws.Cells.Item(11, 3).Value = "NSE:DHUNINV"

# Row 12
$ws.Cells.Item(12, 2).ClearContents()
$ws.Cells.Item(12, 3).Value = "NSE:DREAMFOLKS"

# Row 13
$ws.Cells.Item(13, 2).ClearContents()
$ws.Cells.Item(13, 3).Value = "NSE:ELIN"

# Row 14
$ws.Cells.Item(14, 2).ClearContents()
$ws.Cells.Item(14, 3).Value = "NSE:EVEREADY"

# Row 15
$ws.Cells.Item(15, 2).ClearContents()
$ws.Cells.Item(15, 3).Value = "NSE:GESHIP"

# Row 16
$ws.Cells.Item(16, 2).ClearContents()
$ws.Cells.Item(16, 3).Value = "NSE:GOCLCORP"

# Row 17
$ws.Cells.Item(17, 2).ClearContents()
$ws.Cells.Item(17, 3).Value = "NSE:GPPL"

# Row 18
$ws.Cells.Item(18, 2).ClearContents()
$ws.Cells.Item(18, 3).Value = "NSE:GROBTEA"

# Row 19
$ws.Cells.Item(19, 2).ClearContents()
$ws.Cells.Item(19, 3).Value = "NSE:HINDWAREAP"

# Row 20
$ws.Cells.Item(20, 2).ClearContents()
$ws.Cells.Item(20, 3).Value = "NSE:HLEGLAS"

# Row 21
$ws.Cells.Item(21, 2).ClearContents()
$ws.Cells.Item(21, 3).Value = "NSE:IIFL"

# Row 22
$ws.Cells.Item(22, 2).ClearContents()
$ws.Cells.Item(22, 3).Value = "NSE:INDNIPPON"

# Row 23
$ws.Cells.Item(23, 2).ClearContents()
$ws.Cells.Item(23, 3).Value = "NSE:INDOSTAR"

# Row 24
$ws.Cells.Item(24, 1).Value = 22
$ws.Cells.Item(24, 3).Value = "NSE:IZMO"

# Row 25
$ws.Cells.Item(25, 1).Value = 23
$ws.Cells.Item(25, 3).Value = "NSE:JASH"

# Row 26
$ws.Cells.Item(26, 1).Value = 24
$ws.Cells.Item(26, 3).Value = "NSE:KHAICHEM"

# Row 27
$ws.Cells.Item(27, 1).Value = 25
$ws.Cells.Item(27, 3).Value = "NSE:KNRCON"

# Row 28
$ws.Cells.Item(28, 1).Value = 26
$ws.Cells.Item(28, 3).Value = "NSE:MALLCOM"

# Row 29
$ws.Cells.Item(29, 1).Value = 27
$ws.Cells.Item(29, 3).Value = "NSE:NELCO"

# Row 30
$ws.Cells.Item(30, 1).Value = 28
$ws.Cells.Item(30, 3).Value = "NSE:PARSVNATH"

# Row 31
$ws.Cells.Item(31, 1).Value = 29
$ws.Cells.Item(31, 3).Value = "NSE:PGEL"

# Row 32
$ws.Cells.Item(32, 1).Value = 30
$ws.Cells.Item(32, 3).Value = "NSE:REMSONSIND"

# Row 33
$ws.Cells.Item(33, 1).Value = 31
$ws.Cells.Item(33, 3).Value = "NSE:ROLEXRINGS"

# Row 34
$ws.Cells.Item(34, 1).Value = 32
$ws.Cells.Item(34, 3).Value = "NSE:SAKHTISUG"
